$wb = $excel.ActiveWorkbook

# --- Sheet: u_MAB ---
$ws1 = $wb.Worksheets.Item("u_MAB")

$ws1.Range("B15").Value = 0.06536809121972878

$ws1.Range("A16").Value = 0.05412794201161027
$ws1.Range("B16").Value = 0.07830710215390414

$ws1.Range("A24").Value = 0.3677267441372666
$ws1.Range("B24").Value = 0.08832357937858312

$ws1.Range("B25").Value = 0.2340227058923388

$ws1.Range("A27").Value = 0.1183802703740772
$ws1.Range("B27").Value = 0.08016175730869057

$ws1.Range("A40").Value = 0
$ws1.Range("B40").Value = 0

$ws1.Range("A47").Value = 0.140767149045204
$ws1.Range("B47").Value = 1.746325518992349

$ws1.Range("A49").Value = 0.157809606073953

$ws1.Range("B51").Value = 0.09625705141564801

$ws1.Range("A52").Value = 0.05182702263477508

$ws1.Range("B61").Value = 0

# --- Sheet: u_EOH ---
$ws2 = $wb.Worksheets.Item("u_EOH")

$ws2.Range("A2").Value = -0.2883272966876143
$ws2.Range("A3").Value = -0.1957672534221044

# --- Sheet: v_l ---
$ws3 = $wb.Worksheets.Item("v_l")

$ws3.Range("A2").Value = 3099722.89946284
$ws3.Range("A3").Value = 2111876.148112711
$ws3.Range("A4").Value = 0
